$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 318 (pushes old rows 318-342 down to 319-343).
$ws.Rows.Item(318).Insert()

# The insert carries formatting over from the row above for columns C and V;
# V is not used by the new row, so clear it out to avoid a stray styled cell.
$ws.Range("V318").ClearFormats()
$ws.Range("V318").ClearContents()

# Fill in the new "UK Election Statistics" row.
$ws.Range("A318").Value = "UK Election Statistics"
$ws.Range("B318").Value = "elections"
$ws.Range("C318").Value = "https://researchbriefings.parliament.uk/ResearchBriefing/Summary/CBP-7529"
$ws.Range("D318").Value = "elections, turnout"
$ws.Range("E318").Value = "UK"
$ws.Range("F318").Value = 0
$ws.Range("G318").Value = 0
$ws.Range("H318").Value = 0
$ws.Range("I318").Value = 0
$ws.Range("J318").Value = 1
$ws.Range("K318").Value = 1918
$ws.Range("L318").Value = 2019
$ws.Range("M318").Value = "online"
$ws.Range("N318").Value = "no"
$ws.Range("O318").Value = 1
$ws.Range("P318").Value = "http://researchbriefings.files.parliament.uk/documents/CBP-7529/CBP-7529.pdf"
$ws.Range("T318").Value = "http://researchbriefings.files.parliament.uk/documents/CBP-7529/CBP-7529.Download.xlsx"
$ws.Range("AB318").Value = 20190819

# Add hyperlinks (link / file_codebook / file_excel columns), matching the
# "Link" style already used throughout the sheet for hyperlinked cells.
$ws.Hyperlinks.Add($ws.Range("C318"), "https://researchbriefings.parliament.uk/ResearchBriefing/Summary/CBP-7529")
$ws.Range("C318").Style = "Link"

$ws.Hyperlinks.Add($ws.Range("P318"), "http://researchbriefings.files.parliament.uk/documents/CBP-7529/CBP-7529.pdf")
$ws.Range("P318").Style = "Link"

$ws.Hyperlinks.Add($ws.Range("T318"), "http://researchbriefings.files.parliament.uk/documents/CBP-7529/CBP-7529.Download.xlsx")
$ws.Range("T318").Style = "Link"

# Restore the view state (scrolled/selected near the bottom of the list).
$ws.Application.ActiveWindow.ScrollRow = 307
[void]$ws.Range("J318").Select()
